$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{Row=2;  A=46563; B="Thales Lima"; C="Vendas"; D="Outros"; E=4; F=45085; G=2258.68},
    @{Row=3;  A=27051; B="Ana Clara Silveira"; C="Vendas"; D="Viagem de negocios"; E=2; F=45098; G=8385.379999999999},
    @{Row=4;  A=36088; B="Leonardo Fonseca"; C="Marketing"; D="Consulta medica"; E=8; F=45098; G=9872.65},
    @{Row=5;  A=26627; B="Anthony Lopes"; C="TI"; D="Consulta medica"; E=3; F=45101; G=9464.059999999999},
    @{Row=6;  A=45604; B="Esther Macedo"; C="Atendimento ao Cliente"; D="Consulta medica"; E=3; F=45101; G=6244.86},
    @{Row=7;  A=86661; B="Rafael Sousa"; C="Marketing"; D="Problemas pessoais"; E=2; F=45088; G=4168.05},
    @{Row=8;  A=53867; B="Ana Beatriz Cavalcante"; C="Juridico"; D="Consulta medica"; E=1; F=45092; G=2760.88},
    @{Row=9;  A=72573; B="Cauê Siqueira"; C="Marketing"; D="Viagem de negocios"; E=2; F=45088; G=7907.38},
    @{Row=10; A=77588; B="Benicio Lima"; C="P&D"; D="Viagem de negocios"; E=2; F=45086; G=5915.16},
    @{Row=11; A=37102; B="Bernardo das Neves"; C="Atendimento ao Cliente"; D="Doenca"; E=3; F=45099; G=2496.34}
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 1).Value = $item.A
    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 3).Value = $item.C
    $ws.Cells.Item($r, 4).Value = $item.D
    $ws.Cells.Item($r, 5).Value = $item.E
    $ws.Cells.Item($r, 6).Value = $item.F
    $ws.Cells.Item($r, 7).Value = $item.G
}
